$d = $word.ActiveDocument

$d.Content.Find.Execute("853×8=6824", $true, $true, $false, $false, $false, $true, 1, $false, "711×4=2844", 2) | Out-Null
$d.Content.Find.Execute("638×7=4466", $true, $true, $false, $false, $false, $true, 1, $false, "461×8=3688", 2) | Out-Null
$d.Content.Find.Execute("364×8=2912", $true, $true, $false, $false, $false, $true, 1, $false, "495×6=2970", 2) | Out-Null
$d.Content.Find.Execute("726×3=2178", $true, $true, $false, $false, $false, $true, 1, $false, "493×3=1479", 2) | Out-Null
$d.Content.Find.Execute("418×6=2508", $true, $true, $false, $false, $false, $true, 1, $false, "284×9=2556", 2) | Out-Null
$d.Content.Find.Execute("864×4=3456", $true, $true, $false, $false, $false, $true, 1, $false, "999×4=3996", 2) | Out-Null
$d.Content.Find.Execute("508×5=2540", $true, $true, $false, $false, $false, $true, 1, $false, "746×5=3730", 2) | Out-Null
$d.Content.Find.Execute("258×5=1290", $true, $true, $false, $false, $false, $true, 1, $false, "364×3=1092", 2) | Out-Null
$d.Content.Find.Execute("975×5=4875", $true, $true, $false, $false, $false, $true, 1, $false, "609×2=1218", 2) | Out-Null
$d.Content.Find.Execute("570×5=2850", $true, $true, $false, $false, $false, $true, 1, $false, "120×4=480", 2) | Out-Null
$d.Content.Find.Execute("655×9=5895", $true, $true, $false, $false, $false, $true, 1, $false, "577×4=2308", 2) | Out-Null
$d.Content.Find.Execute("615×5=3075", $true, $true, $false, $false, $false, $true, 1, $false, "113×8=904", 2) | Out-Null
$d.Content.Find.Execute("219×8=1752", $true, $true, $false, $false, $false, $true, 1, $false, "767×7=5369", 2) | Out-Null
$d.Content.Find.Execute("289×6=1734", $true, $true, $false, $false, $false, $true, 1, $false, "166×6=996", 2) | Out-Null
$d.Content.Find.Execute("777×8=6216", $true, $true, $false, $false, $false, $true, 1, $false, "470×4=1880", 2) | Out-Null
$d.Content.Find.Execute("684×9=6156", $true, $true, $false, $false, $false, $true, 1, $false, "630×6=3780", 2) | Out-Null
$d.Content.Find.Execute("773×8=6184", $true, $true, $false, $false, $false, $true, 1, $false, "511×4=2044", 2) | Out-Null
$d.Content.Find.Execute("121×9=1089", $true, $true, $false, $false, $false, $true, 1, $false, "154×6=924", 2) | Out-Null
$d.Content.Find.Execute("756×8=6048", $true, $true, $false, $false, $false, $true, 1, $false, "816×2=1632", 2) | Out-Null
$d.Content.Find.Execute("435×2=870", $true, $true, $false, $false, $false, $true, 1, $false, "517×7=3619", 2) | Out-Null
$d.Content.Find.Execute("809×2=1618", $true, $true, $false, $false, $false, $true, 1, $false, "327×7=2289", 2) | Out-Null
$d.Content.Find.Execute("436×4=1744", $true, $true, $false, $false, $false, $true, 1, $false, "937×8=7496", 2) | Out-Null
$d.Content.Find.Execute("525×2=1050", $true, $true, $false, $false, $false, $true, 1, $false, "402×3=1206", 2) | Out-Null
$d.Content.Find.Execute("977×6=5862", $true, $true, $false, $false, $false, $true, 1, $false, "552×8=4416", 2) | Out-Null
$d.Content.Find.Execute("468×4=1872", $true, $true, $false, $false, $false, $true, 1, $false, "679×6=4074", 2) | Out-Null
